$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("General_&_Special_Internal_1")

# Set column G (Duration) to 900 for every data row (rows 2 through 276)
$ws.Range("G2:G276").Value = 900
